$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.552.53"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "2.717.09"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  +0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "559.34"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.81%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "157.42"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  +0.00%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.592"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("E10").Value = "  -0.39%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.63"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.72%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.372"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.99%  "

$ws.Range("D13").Value = "3.197.40"
$ws.Range("E13").Value = "  -1.46%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.49"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "

$ws.Range("D15").Value = "63.401.14"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").Value = "2.718.41"
$ws.Range("E17").Value = "  -1.56%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "12.17"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.65"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.36%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "349.50"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "

$ws.Range("E21").Value = "  -4.36%  "

$ws.Range("E22").Value = "  +0.07%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.514"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.86%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "64.17"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.170"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").Value = "  +0.03%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.17"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.84%  "

$ws.Range("D28").Value = "0.0₃0879"
$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.37"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +8.85%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.15"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "164.29"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.86"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.45%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.82"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.46"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("E37").Value = "  -0.60%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "348.27"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.958"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.63%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.08"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.09%  "

$ws.Range("E41").Value = "  -4.66%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "38.19"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.23%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "21.29"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.88%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "20.67"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.63%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0571"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.14%  "

$ws.Range("E46").Value = "  -1.43%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "131.87"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.21%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "11.08"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0983"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.29%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0245"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.77%  "
